$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "coldata" to "Folha1"
$ws.Name = "Folha1"

# Remove the last 10 data rows (the "C" group, rows 22-31), shifting rows up
$ws.Range("A22:D31").EntireRow.Delete()

# Match the author's final selection in the sheet
$ws.Range("D22").Select()
